$d = $word.ActiveDocument

# Locate the end of the "Historia de usuario no.2" body paragraph
# ("...sin tener que eliminarlos por completo.") using Find, then expand
# to the full paragraph so we get a reliable anchor regardless of exact
# character offsets.
$rng = $d.Content
$found = $rng.Find.Execute("sin tener que eliminarlos por completo.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text"
}
[void]$rng.Collapse(0)
[void]$rng.Expand(4)
$insertAt = $rng.End

$target = $d.Range($insertAt, $insertAt)

# Build the OOXML fragment to insert. It contains (in order):
#   - 3 blank "Sin espaciado" paragraphs (matching the style already used
#     for the blank separator paragraph that sits at the end of the doc)
#   - the new "Historia de Usuario no.3" title paragraph
#   - the new body paragraph (4 runs, incl. a manual line break)
#   - a trailing throw-away empty paragraph
#
# NOTE: when OOXML is inserted at a collapsed Range via InsertXML, every
# paragraph mark in the fragment becomes a genuine new paragraph EXCEPT
# the very last one, whose run content (if any) gets folded into the
# paragraph that already began at the insertion point (keeping that
# paragraph's own pPr). Appending a throw-away empty <w:p/> as the last
# element means that "merge" consumes the throw-away instead of our real
# content, so every real paragraph above it keeps its own exact pPr. The
# throw-away paragraph (now a real, but unwanted, empty paragraph) is
# deleted afterwards.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="5670"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="5670"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="Sinespaciado"/><w:tabs><w:tab w:val="left" w:pos="5670"/></w:tabs><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>Historia de Usuario no.3: Consultar seguros vigentes por cliente o por vehículo.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Como el asesor del área de seguros, quiero poder consultar fácilmente todos los seguros vigentes que se han registrado en el sistema, para verificar rápidamente la cobertura, vigencia, aseguradora y número de póliza, </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>seguradora, tipo de seguro y estado actual.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/><w:t>El sistema debe permitirme filtrar por estado del seguro (activo, vencido), por vehículo (placa), y por cliente si es necesario, mostrando los detalles relevantes para auditorías, renovaciones o reportes internos.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)

# Remove the throw-away empty paragraph: it's the paragraph that now
# immediately follows the body paragraph we just inserted (found again
# via its distinctive closing text) and precedes the document's
# pre-existing trailing blank paragraph.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("renovaciones o reportes internos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find newly inserted body paragraph"
}
[void]$rng2.Collapse(0)
[void]$rng2.Expand(4)
$afterBody = $d.Range($rng2.End, $rng2.End)
[void]$afterBody.Expand(4)
[void]$afterBody.Delete()
